# Commit: swap the presentation's table style to a theme-based style,
# and swap the Integral/Office Theme colour schemes between the deck's
# two theme parts (theme1.xml <-> theme2.xml).
#
# theme1.xml currently carries the "Integral" / "Red Violet" palette and
# is the theme actually used by the slide master (and therefore every
# slide). theme2.xml carries the stock "Office Theme" / "Office" palette
# and is only used by the notes master. The authored edit swaps the two
# files' contents wholesale (names + 12 scheme colours); the fonts and
# format scheme are identical between the two themes already.
#
# The PowerPoint object model only exposes a single writable
# ThemeColorScheme (reached from the slide master / any slide), which is
# backed by theme1.xml, so that is what we repaint here with the target
# "Office" colour values.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Three tables (slides 14, 15, 16) move off the deck's custom
#    "Table_0" style onto the built-in theme table style.
# ---------------------------------------------------------------------
$newTableStyleId = "{009BA083-EC07-41CA-A91D-89F46884716D}"

$tableSlideIndexes = @(14, 15, 16)
foreach ($slideIndex in $tableSlideIndexes) {
    $slide = $p.Slides.Item($slideIndex)
    foreach ($shape in $slide.Shapes) {
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($newTableStyleId)
        }
    }
}

# ---------------------------------------------------------------------
# 2) Repaint the master theme's 12 scheme colours from the "Integral" /
#    Red Violet palette to the stock "Office Theme" / Office palette
#    (this is the palette that theme2.xml already carries, and what
#    theme1.xml ends up with after the swap).
#    Index order: dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink.
#    COM RGB values are packed 0x00BBGGRR (i.e. decimal of BBGGRR).
# ---------------------------------------------------------------------
$master = $p.SlideMaster
$theme = $master.Theme
$colors = $theme.ThemeColorScheme

$officePalette = @(
    0,          # dk1      000000
    16777215,   # lt1      FFFFFF
    6968388,    # dk2      44546A
    15132391,   # lt2      E7E6E6
    13998939,   # accent1  5B9BD5
    3243501,    # accent2  ED7D31
    10855845,   # accent3  A5A5A5
    49407,      # accent4  FFC000
    12874308,   # accent5  4472C4
    4697456,    # accent6  70AD47
    12673797,   # hlink    0563C1
    7491477     # folHlink 954F72
)

for ($i = 1; $i -le $colors.Count; $i++) {
    $colors.Item($i).RGB = $officePalette[$i - 1]
}
